$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1217.2174
$ws.Range("J17").Value = 998.42224
$ws.Range("L17").Value = 2995.26672
$ws.Range("N17").Value = -3331.26672

$ws.Range("H101").Value = 1699.5
$ws.Range("I101").Value = 400
$ws.Range("J101").Value = 2999
$ws.Range("K101").Value = 1200
$ws.Range("L101").Value = 8997
$ws.Range("M101").Value = 422
$ws.Range("N101").Value = -12241

$ws.Range("H137").Value = 1221.25
$ws.Range("I137").Value = 1167.1428
$ws.Range("K137").Value = 3501.4284
$ws.Range("M137").Value = -951.4284000000002

$ws.Range("H138").Value = 2812.5881
$ws.Range("I138").Value = 2768.52
$ws.Range("J138").Value = 2854.9614
$ws.Range("K138").Value = 8305.559999999999
$ws.Range("L138").Value = 8564.8842
$ws.Range("M138").Value = -3165.559999999999
$ws.Range("N138").Value = -18844.8842

$ws.Range("H141").Value = 2004035.6
$ws.Range("I141").Value = 3114268.2
$ws.Range("J141").Value = 5617
$ws.Range("K141").Value = 9342804.600000001
$ws.Range("L141").Value = 16851
$ws.Range("M141").Value = -9337624.600000001
$ws.Range("N141").Value = -27211

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 148.75
$ws.Range("I5").Value = 148.75
$ws.Range("K5").Value = 148.75
$ws.Range("M5").Value = -36.75

$ws.Range("H61").Value = 3852.4062
$ws.Range("I61").Value = 2683.2856
$ws.Range("J61").Value = 12036.25
$ws.Range("K61").Value = 2683.2856
$ws.Range("L61").Value = 12036.25
$ws.Range("M61").Value = -2471.2856
$ws.Range("N61").Value = -12460.25

$ws.Range("H122").Value = 7748.25
$ws.Range("I122").Value = 7748.25
$ws.Range("K122").Value = 23244.75
$ws.Range("M122").Value = -20794.75

$ws.Range("H136").Value = 3852.4062
$ws.Range("I136").Value = 2683.2856
$ws.Range("J136").Value = 12036.25
$ws.Range("K136").Value = 8049.8568
$ws.Range("L136").Value = 36108.75
$ws.Range("M136").Value = -5499.8568
$ws.Range("N136").Value = -41208.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 148.75
$ws.Range("I4").Value = 148.75
$ws.Range("K4").Value = 148.75
$ws.Range("M4").Value = -33.75

$ws.Range("H86").Value = 107487.63
$ws.Range("J86").Value = 335716.5
$ws.Range("L86").Value = 335716.5
$ws.Range("N86").Value = -337962.5

$ws.Range("H89").Value = 107487.63
$ws.Range("J89").Value = 335716.5
$ws.Range("L89").Value = 1678582.5
$ws.Range("N89").Value = -1689814.5

$ws.Range("H99").Value = 1515.5
$ws.Range("I99").Value = 1242.8889
$ws.Range("K99").Value = 1242.8889
$ws.Range("M99").Value = 255.1111000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2828.639
$ws.Range("I31").Value = 1902.6086
$ws.Range("K31").Value = 1902.6086
$ws.Range("M31").Value = -1607.6086

$ws.Range("H34").Value = 2828.639
$ws.Range("I34").Value = 1902.6086
$ws.Range("K34").Value = 1902.6086
$ws.Range("M34").Value = -1700.6086

$ws.Range("H99").Value = 1252614
$ws.Range("J99").Value = 3152.1667
$ws.Range("L99").Value = 3152.1667
$ws.Range("N99").Value = -6148.1667

$ws.Range("H126").Value = 1252614
$ws.Range("J126").Value = 3152.1667
$ws.Range("L126").Value = 9456.500100000001
$ws.Range("N126").Value = -14396.5001

$ws.Range("H132").Value = 2628
$ws.Range("I132").Value = 1250.375
$ws.Range("J132").Value = 3730.1
$ws.Range("K132").Value = 3751.125
$ws.Range("L132").Value = 11190.3
$ws.Range("M132").Value = -1221.125
$ws.Range("N132").Value = -16250.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 700
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H107").Value = 906.3333
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 906.3333
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").Value = 2718.9999
$ws.Range("N107").Value = -6558.9999

$ws.Range("H112").Value = 4030
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 4030
$ws.Range("K112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("M112").Value = 12090
$ws.Range("N112").Value = -14306

$ws.Range("H130").Value = 2083.3333
$ws.Range("J130").Value = 2083.3333
$ws.Range("L130").Value = 6249.999899999999
$ws.Range("N130").Value = -16289.9999

$ws.Range("H131").Value = 11188.299
$ws.Range("J131").Value = 12782.403
$ws.Range("L131").Value = 38347.209
$ws.Range("N131").Value = -48427.209

$ws.Range("H136").Value = 3584.6
$ws.Range("I136").Value = 3584.6
$ws.Range("K136").Value = 10753.8
$ws.Range("M136").Value = -5653.799999999999

$ws.Range("H140").Value = 1676.7142
$ws.Range("I140").Value = 977.0769
$ws.Range("J140").Value = 2283.0667
$ws.Range("K140").Value = 2931.2307
$ws.Range("L140").Value = 6849.2001
$ws.Range("M140").Value = 2248.7693
$ws.Range("N140").Value = -17209.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("N26").Value = 0

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("N50").Value = 0

$ws.Range("H102").Value = 2664.7334
$ws.Range("I102").Value = 2650.5386
$ws.Range("K102").Value = 2650.5386
$ws.Range("M102").Value = -1028.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 245173.08
$ws.Range("J2").Value = 26750
$ws.Range("L2").Value = 26750
$ws.Range("N2").Value = -26974

$ws.Range("H40").Value = 5755.5884
$ws.Range("I40").Value = 1384.8
$ws.Range("J40").Value = 11999.571
$ws.Range("K40").Value = 1384.8
$ws.Range("L40").Value = 11999.571
$ws.Range("M40").Value = -1248.8
$ws.Range("N40").Value = -12271.571

$ws.Range("H122").Value = 7405.273
$ws.Range("I122").Value = 5181.625
$ws.Range("J122").Value = 13335
$ws.Range("K122").Value = 15544.875
$ws.Range("L122").Value = 40005
$ws.Range("M122").Value = -13094.875
$ws.Range("N122").Value = -44905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 606.92
$ws.Range("I107").Value = 414.05264
$ws.Range("J107").Value = 1217.6666
$ws.Range("K107").Value = 1242.15792
$ws.Range("L107").Value = 3652.9998
$ws.Range("M107").Value = 677.8420799999999
$ws.Range("N107").Value = -7492.9998

$ws.Range("H122").Value = 214528.33
$ws.Range("I122").Value = 257034
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 771102
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -768652
$ws.Range("N122").Value = -10900

$ws.Range("H126").Value = 20538
$ws.Range("I126").Value = 26860.8
$ws.Range("K126").Value = 80582.39999999999
$ws.Range("M126").Value = -78112.39999999999
